$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "NEW High Ambition XML Name" / "NEW Low Ambition XML Name" values
# for the Uncertainty_2 (EV Costs) row - F14:G14
$ws.Range("F14:G14").ClearContents()

# Rename the Uncertainty_6 "NEW High Ambition XML Name" value (E18)
$ws.Range("E18").Value = "Uncertainty_6_High_HOV_sw_high_CL"

# Update the saved selection/active cell to F14
$ws.Range("F14").Select()
